# Swap the contents of row 2 and row 3 for every column that is actually
# used by either row (A, B, D-I, K, L, M, N, P-W, Y, AA, AC-AE, AG, AT,
# AW-AY). Columns that are not used by either row are left completely
# untouched so no stray blank cells are introduced.
#
# A scratch cell far outside the used range (row 500) is used as a
# temporary holder so that values can be exchanged safely, and
# Range.Copy is used (instead of re-typing values) so that cell types
# (numbers, booleans, text that looks like a date, etc.) are preserved
# exactly instead of being re-interpreted by Excel's input parser.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratchRow = 500

# Columns present in row 2 (before the edit).
$row2Cols = @("A","B","D","E","F","G","H","I","K","L","M","N","P","Q","R","S","T","U","V","W","Y","AA","AC","AD","AE","AG","AT","AW","AX","AY")
# Columns present in row 3 (before the edit).
$row3Cols = @("A","B","D","E","F","G","H","I","M","P","Q","R","S","T","U","V","W","Y","AA","AC","AD","AE","AG","AT","AW","AX","AY")

foreach ($col in $row2Cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $scratch = $ws.Range($col + $scratchRow)

    if ($row3Cols -contains $col) {
        # Cell exists in both rows: swap the two values via the scratch cell.
        $cell2.Copy($scratch)
        $cell3.Copy($cell2)
        $scratch.Copy($cell3)
        $scratch.Clear()
    } else {
        # Cell exists only in row 2: move it to row 3 and remove it from row 2.
        $cell2.Copy($cell3)
        $cell2.Clear()
    }
}
